{"js": "// Replace the date and the ten three-digit-by-one-digit multiplication\n// expressions with their updated values, exactly as described by the diff.\nconst replacements = [\n  [\"2025-08-16 Saturday\", \"2025-08-17 Sunday\"],\n  [\"248\u00d72=\", \"599\u00d74=\"],\n  [\"644\u00d76=\", \"562\u00d76=\"],\n  [\"876\u00d72=\", \"506\u00d78=\"],\n  [\"637\u00d75=\", \"588\u00d79=\"],\n  [\"883\u00d74=\", \"135\u00d73=\"],\n  [\"719\u00d75=\", \"268\u00d74=\"],\n  [\"451\u00d74=\", \"292\u00d76=\"],\n  [\"340\u00d79=\", \"933\u00d76=\"],\n  [\"308\u00d76=\", \"240\u00d79=\"],\n  [\"123\u00d73=\", \"901\u00d72=\"],\n  [\"954\u00d78=\", \"593\u00d76=\"],\n  [\"929\u00d77=\", \"510\u00d76=\"],\n  [\"440\u00d78=\", \"935\u00d74=\"],\n  [\"568\u00d73=\", \"684\u00d72=\"],\n  [\"602\u00d74=\", \"907\u00d72=\"],\n  [\"554\u00d73=\", \"853\u00d79=\"],\n  [\"308\u00d77=\", \"514\u00d79=\"],\n  [\"248\u00d78=\", \"545\u00d79=\"],\n  [\"360\u00d72=\", \"243\u00d74=\"],\n  [\"370\u00d75=\", \"941\u00d77=\"],\n  [\"919\u00d79=\", \"745\u00d73=\"],\n  [\"271\u00d78=\", \"444\u00d74=\"],\n  [\"756\u00d73=\", \"606\u00d72=\"],\n  [\"771\u00d76=\", \"829\u00d76=\"],\n  [\"204\u00d77=\", \"801\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the ten three-digit-by-one-digit multiplication\n# expressions with their new values, exactly as described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-16 Saturday\", \"2025-08-17 Sunday\"),\n    @(\"248\u00d72=\", \"599\u00d74=\"),\n    @(\"644\u00d76=\", \"562\u00d76=\"),\n    @(\"876\u00d72=\", \"506\u00d78=\"),\n    @(\"637\u00d75=\", \"588\u00d79=\"),\n    @(\"883\u00d74=\", \"135\u00d73=\"),\n    @(\"719\u00d75=\", \"268\u00d74=\"),\n    @(\"451\u00d74=\", \"292\u00d76=\"),\n    @(\"340\u00d79=\", \"933\u00d76=\"),\n    @(\"308\u00d76=\", \"240\u00d79=\"),\n    @(\"123\u00d73=\", \"901\u00d72=\"),\n    @(\"954\u00d78=\", \"593\u00d76=\"),\n    @(\"929\u00d77=\", \"510\u00d76=\"),\n    @(\"440\u00d78=\", \"935\u00d74=\"),\n    @(\"568\u00d73=\", \"684\u00d72=\"),\n    @(\"602\u00d74=\", \"907\u00d72=\"),\n    @(\"554\u00d73=\", \"853\u00d79=\"),\n    @(\"308\u00d77=\", \"514\u00d79=\"),\n    @(\"248\u00d78=\", \"545\u00d79=\"),\n    @(\"360\u00d72=\", \"243\u00d74=\"),\n    @(\"370\u00d75=\", \"941\u00d77=\"),\n    @(\"919\u00d79=\", \"745\u00d73=\"),\n    @(\"271\u00d78=\", \"444\u00d74=\"),\n    @(\"756\u00d73=\", \"606\u00d72=\"),\n    @(\"771\u00d76=\", \"829\u00d76=\"),\n    @(\"204\u00d77=\", \"801\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
